$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Keep F1 blank (it is an empty placeholder cell in the original sheet;
# re-assert emptiness so the round-trip save doesn't coerce it to a value)
$ws.Cells.Item(1, 6).Value = ""

# Add "NA" under duplicate_image_filename (column E) for data rows 2-21
for ($r = 2; $r -le 21; $r++) {
    $ws.Cells.Item($r, 5).Value = "NA"
}
